$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 68: 2026-02-01 bitcoin buy.
# Column A mirrors the recent rows (61-67): plain text date string, default
# (unstyled) cell. A leading apostrophe forces Excel to store the value as
# literal text instead of auto-converting "02/01/2026" into a date serial;
# resetting the style to "Normal" afterwards drops the quote-prefix
# formatting flag so the cell ends up with no explicit style, matching the
# rest of the recent rows.
$ws.Range("A68").Value = "'02/01/2026"
$ws.Range("A68").Style = "Normal"

$ws.Range("B68").Value = 0.0006264100000000078
$ws.Range("C68").Value = 79021.72698392329
$ws.Range("D68").Value = 50
